# "Fixed spacing in diagram"
# Nudge the positions of four small "Graphic" icon pictures on the
# architecture-diagram slide. The diff's <a:off> values are EMUs
# (914400 EMU = 1 inch); PowerPoint's Shape.Left/Top COM properties
# are expressed in points (72 pt = 1 inch), so EMU / 12700 = points.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Map of shape Id -> new (x,y) offset in EMUs, taken from the target OOXML.
$targets = @{
    74 = @(5051641, 5460193)   # "Graphic 24"
    76 = @(5034336, 3424825)   # "Graphic 8"
    77 = @(5086967, 4185653)   # "Graphic 76"
    78 = @(5091670, 6167816)   # "Graphic 77"
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($targets.ContainsKey($shp.Id)) {
        $coords = $targets[$shp.Id]
        $shp.Left = $coords[0] / 12700
        $shp.Top  = $coords[1] / 12700
    }
}
